$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '23.272.53'
$ws.Range('E2').Value = '  -0.79%  '
$ws.Range('D3').Value = '1.615.15'
$ws.Range('E3').Value = '  -1.11%  '
$ws.Range('D4').Value = "'0.9995"
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = "'0.9996"
$ws.Range('E5').Value = '  -0.11%  '
$ws.Range('D6').Value = "'302.33"
$ws.Range('E6').Value = '  -0.85%  '
$ws.Range('D7').Value = "'0.3739"
$ws.Range('E7').Value = '  -0.18%  '
$ws.Range('D8').Value = "'51.98"
$ws.Range('E8').Value = '  +0.93%  '
$ws.Range('D9').Value = "'0.3540"
$ws.Range('E9').Value = '  -2.68%  '
$ws.Range('D10').Value = "'0.08133"
$ws.Range('E10').Value = '  -0.61%  '
$ws.Range('D11').Value = "'1.212"
$ws.Range('E11').Value = '  -1.16%  '
$ws.Range('D12').Value = "'0.9993"
$ws.Range('E12').Value = '  -0.11%  '
$ws.Range('D13').Value = "'22.04"
$ws.Range('E13').Value = '  -2.38%  '
$ws.Range('D14').Value = "'6.424"
$ws.Range('E14').Value = '  -1.88%  '
$ws.Range('D15').Value = "'7.232"
$ws.Range('E15').Value = '  -0.18%  '
$ws.Range('D16').Value = "'0.00001216"
$ws.Range('E16').Value = '  -2.73%  '
$ws.Range('D17').Value = '1.614.23'
$ws.Range('E17').Value = '  -0.97%  '
$ws.Range('D18').Value = "'94.72"
$ws.Range('E18').Value = '  +0.25%  '
$ws.Range('D19').Value = "'0.06920"
$ws.Range('E19').Value = '  -0.85%  '
$ws.Range('D20').Value = "'6.562"
$ws.Range('E20').Value = '  +1.61%  '
$ws.Range('D21').Value = "'17.29"
$ws.Range('E21').Value = '  -2.68%  '
$ws.Range('D22').Value = "'0.9997"
$ws.Range('E22').Value = '  -0.20%  '
$ws.Range('D23').Value = "'12.38"
$ws.Range('E23').Value = '  -2.87%  '
$ws.Range('B24').Value = 'WrappedBTC'
$ws.Range('C24').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D24').Value = '23.258.98'
$ws.Range('E24').Value = '  -0.83%  '
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').Value = "'2.507"
$ws.Range('E25').Value = '  +1.44%  '
$ws.Range('D26').Value = "'3.078"
$ws.Range('E26').Value = '  -4.10%  '
$ws.Range('D27').Value = "'20.89"
$ws.Range('E27').Value = '  -2.12%  '
$ws.Range('D28').Value = "'152.27"
$ws.Range('E28').Value = '  +1.64%  '
$ws.Range('D29').Value = "'5.241"
$ws.Range('E29').Value = '  -2.34%  '
$ws.Range('D30').Value = "'133.05"
$ws.Range('E30').Value = '  -1.40%  '
$ws.Range('D31').Value = '1.792.65'
$ws.Range('E31').Value = '  -1.00%  '
$ws.Range('E32').Value = '  +6.82%  '
$ws.Range('D33').Value = "'6.508"
$ws.Range('E33').Value = '  -4.65%  '
$ws.Range('B34').Value = 'FraxShare'
$ws.Range('C34').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D34').Value = "'11.64"
$ws.Range('E34').Value = '  +5.07%  '
$ws.Range('B35').Value = 'WEMIXTOKEN'
$ws.Range('C35').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D35').Value = "'2.026"
$ws.Range('E35').Value = '  -11.31%  '
$ws.Range('D36').Value = "'0.02720"
$ws.Range('E36').Value = '  -2.43%  '
$ws.Range('D37').Value = "'0.08730"
$ws.Range('E37').Value = '  -0.47%  '
$ws.Range('D38').Value = "'0.2458"
$ws.Range('E38').Value = '  -2.86%  '
$ws.Range('D39').Value = "'0.06933"
$ws.Range('E39').Value = '  -2.98%  '
$ws.Range('D40').Value = "'5.866"
$ws.Range('E40').Value = '  -3.22%  '
$ws.Range('B41').Value = 'Aptos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D41').Value = "'12.44"
$ws.Range('E41').Value = '  +1.04%  '
$ws.Range('D42').Value = "'0.6885"
$ws.Range('E42').Value = '  -2.66%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').Value = "'1.318"
$ws.Range('E43').Value = '  -2.00%  '
$ws.Range('D44').Value = "'15.54"
$ws.Range('E44').Value = '  -4.05%  '
$ws.Range('D45').Value = "'0.9991"
$ws.Range('E45').Value = '  -0.13%  '
$ws.Range('D46').Value = "'0.6331"
$ws.Range('E46').Value = '  -2.89%  '
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').Value = "'2.255"
$ws.Range('E47').Value = '  -3.30%  '
$ws.Range('B48').Value = 'PancakeSwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D48').Value = "'3.934"
$ws.Range('E48').Value = '  -1.41%  '
$ws.Range('D49').Value = "'0.07871"
$ws.Range('E49').Value = '  -2.17%  '
$ws.Range('D50').Value = "'127.07"
$ws.Range('E50').Value = '  +1.83%  '
$ws.Range('D51').Value = "'1.158"
$ws.Range('E51').Value = '  -3.84%  '
